$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("04-09-2021", "05-09-2021", "06-09-2021", "07-09-2021")

$row = 248
foreach ($d in $dates) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Formula = "=""" + $d + """"
    $aCell.Copy()
    $aCell.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = 17537
    $ws.Cells.Item($row, 3).Value = 1456
    $ws.Cells.Item($row, 4).Value = 521

    $row++
}
